$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet / sheet1.xml) - update "想去人数" (column F) counts
$wsExhibit = $wb.Worksheets.Item(1)
$wsExhibit.Range("F4").Value = 631
$wsExhibit.Range("F6").Value = 9500
$wsExhibit.Range("F7").Value = 855
$wsExhibit.Range("F9").Value = 1206
$wsExhibit.Range("F10").Value = 1181
$wsExhibit.Range("F13").Value = 20
$wsExhibit.Range("F18").Value = 1304

# Sheet "全部类型" (4th sheet / sheet4.xml) - same events, rows offset by one
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F5").Value = 631
$wsAll.Range("F7").Value = 9500
$wsAll.Range("F8").Value = 855
$wsAll.Range("F10").Value = 1206
$wsAll.Range("F11").Value = 1181
$wsAll.Range("F14").Value = 20
$wsAll.Range("F19").Value = 1304
